# edit.ps1
# Applies the "Updated cryptos list" data refresh: updates Price (D) and
# Volume(1h) (E) figures for the crypto rows, and swaps the RocketPoolETH /
# Mantle rows (47 <-> 48) to reflect the new ranking order, per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered list of (cell, new value) pairs taken from the diff.
$cellUpdates = New-Object System.Collections.ArrayList
[void]$cellUpdates.Add(@("D2", "28.874.27"))
[void]$cellUpdates.Add(@("E2", "  -1.97%  "))
[void]$cellUpdates.Add(@("D3", "1.832.46"))
[void]$cellUpdates.Add(@("E3", "  -2.08%  "))
[void]$cellUpdates.Add(@("D4", "0.9997"))
[void]$cellUpdates.Add(@("D5", "245.30"))
[void]$cellUpdates.Add(@("E5", "  +0.64%  "))
[void]$cellUpdates.Add(@("D6", "0.6893"))
[void]$cellUpdates.Add(@("E6", "  -2.34%  "))
[void]$cellUpdates.Add(@("E7", "  -0.11%  "))
[void]$cellUpdates.Add(@("D8", "0.07704"))
[void]$cellUpdates.Add(@("E8", "  -2.86%  "))
[void]$cellUpdates.Add(@("D9", "0.3048"))
[void]$cellUpdates.Add(@("E9", "  -3.02%  "))
[void]$cellUpdates.Add(@("D10", "23.37"))
[void]$cellUpdates.Add(@("E10", "  -4.71%  "))
[void]$cellUpdates.Add(@("D11", "0.07803"))
[void]$cellUpdates.Add(@("E11", "  -0.01%  "))
[void]$cellUpdates.Add(@("D12", "1.838.13"))
[void]$cellUpdates.Add(@("E12", "  -1.73%  "))
[void]$cellUpdates.Add(@("D13", "5.088"))
[void]$cellUpdates.Add(@("E13", "  -1.51%  "))
[void]$cellUpdates.Add(@("D14", "90.92"))
[void]$cellUpdates.Add(@("E14", "  -2.93%  "))
[void]$cellUpdates.Add(@("D15", "0.6810"))
[void]$cellUpdates.Add(@("E15", "  -3.02%  "))
[void]$cellUpdates.Add(@("D16", "6.437"))
[void]$cellUpdates.Add(@("E16", "  -0.96%  "))
[void]$cellUpdates.Add(@("D17", "0.000008314"))
[void]$cellUpdates.Add(@("E17", "  -2.69%  "))
[void]$cellUpdates.Add(@("D18", "28.858.80"))
[void]$cellUpdates.Add(@("E18", "  -2.08%  "))
[void]$cellUpdates.Add(@("D19", "242.15"))
[void]$cellUpdates.Add(@("E19", "  -4.45%  "))
[void]$cellUpdates.Add(@("D20", "2.075.88"))
[void]$cellUpdates.Add(@("E20", "  -2.42%  "))
[void]$cellUpdates.Add(@("E21", "  -3.08%  "))
[void]$cellUpdates.Add(@("E22", "  +0.00%  "))
[void]$cellUpdates.Add(@("D23", "7.454"))
[void]$cellUpdates.Add(@("E23", "  -1.99%  "))
[void]$cellUpdates.Add(@("D24", "1.000"))
[void]$cellUpdates.Add(@("D25", "0.1481"))
[void]$cellUpdates.Add(@("E25", "  -3.75%  "))
[void]$cellUpdates.Add(@("D26", "158.46"))
[void]$cellUpdates.Add(@("E26", "  -1.77%  "))
[void]$cellUpdates.Add(@("D27", "8.791"))
[void]$cellUpdates.Add(@("E27", "  -2.43%  "))
[void]$cellUpdates.Add(@("D28", "18.21"))
[void]$cellUpdates.Add(@("E28", "  -2.95%  "))
[void]$cellUpdates.Add(@("D29", "1.544"))
[void]$cellUpdates.Add(@("E29", "  -0.94%  "))
[void]$cellUpdates.Add(@("D30", "4.218"))
[void]$cellUpdates.Add(@("E30", "  -2.09%  "))
[void]$cellUpdates.Add(@("D31", "4.148"))
[void]$cellUpdates.Add(@("E31", "  -2.66%  "))
[void]$cellUpdates.Add(@("D32", "1.192"))
[void]$cellUpdates.Add(@("E32", "  -0.98%  "))
[void]$cellUpdates.Add(@("D33", "0.05101"))
[void]$cellUpdates.Add(@("E33", "  -3.42%  "))
[void]$cellUpdates.Add(@("D34", "0.7769"))
[void]$cellUpdates.Add(@("E34", "  +2.20%  "))
[void]$cellUpdates.Add(@("D35", "1.850"))
[void]$cellUpdates.Add(@("E35", "  -2.45%  "))
[void]$cellUpdates.Add(@("D36", "1.140"))
[void]$cellUpdates.Add(@("E36", "  -3.61%  "))
[void]$cellUpdates.Add(@("D37", "2.693"))
[void]$cellUpdates.Add(@("E37", "  -0.48%  "))
[void]$cellUpdates.Add(@("D38", "0.01850"))
[void]$cellUpdates.Add(@("E38", "  -1.40%  "))
[void]$cellUpdates.Add(@("D39", "1.218.50"))
[void]$cellUpdates.Add(@("E39", "  -4.77%  "))
[void]$cellUpdates.Add(@("D40", "2.694"))
[void]$cellUpdates.Add(@("E40", "  -2.32%  "))
[void]$cellUpdates.Add(@("D41", "0.9573"))
[void]$cellUpdates.Add(@("E41", "  +6.53%  "))
[void]$cellUpdates.Add(@("D42", "108.64"))
[void]$cellUpdates.Add(@("E42", "  -1.03%  "))
[void]$cellUpdates.Add(@("D43", "5.829"))
[void]$cellUpdates.Add(@("E43", "  -2.51%  "))
[void]$cellUpdates.Add(@("D44", "0.9992"))
[void]$cellUpdates.Add(@("E44", "  -0.13%  "))
[void]$cellUpdates.Add(@("D45", "9.630"))
[void]$cellUpdates.Add(@("E45", "  -0.08%  "))
[void]$cellUpdates.Add(@("E46", "  -3.54%  "))
[void]$cellUpdates.Add(@("B47", "Mantle"))
[void]$cellUpdates.Add(@("C47", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"))
[void]$cellUpdates.Add(@("D47", "0.5159"))
[void]$cellUpdates.Add(@("E47", "  -0.23%  "))
[void]$cellUpdates.Add(@("B48", "RocketPoolETH"))
[void]$cellUpdates.Add(@("C48", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"))
[void]$cellUpdates.Add(@("D48", "1.975.33"))
[void]$cellUpdates.Add(@("E48", "  -2.80%  "))
[void]$cellUpdates.Add(@("D49", "64.38"))
[void]$cellUpdates.Add(@("E49", "  -9.14%  "))
[void]$cellUpdates.Add(@("D50", "1.746"))
[void]$cellUpdates.Add(@("E50", "  -3.12%  "))
[void]$cellUpdates.Add(@("D51", "0.05910"))
[void]$cellUpdates.Add(@("E51", "  -1.31%  "))

foreach ($pair in $cellUpdates) {
    $cellRef = $pair[0]
    $newValue = $pair[1]
    $range = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "1.000", "0.9997",
    # "28.874.27") round-trip as text instead of being coerced to numbers,
    # matching how these cells are stored in the source workbook.
    $range.NumberFormat = "@"
    $range.Value = $newValue
}
